# Sprint 46 - Day 9 Test Case Summary: fill in the totals now that the
# test cases for the Android purchase page UI have been written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

$ws.Range("C51").Value = 7322
$ws.Range("C52").Value = 3220
$ws.Range("C53").Value = 3220

$ws.Range("C53").Select()
